# Fruta / hortaliza, semanal
# A new weekly record is inserted at row 81 ("Cultivar IV Región" / Chirimoya
# data for date 44523), and every existing record from row 81 onward shifts
# down by one row (old row 81 -> new row 82, ..., old row 111 -> new row 112).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 81; this pushes rows 81..111 down to 82..112
$ws.Rows.Item(81).Insert()

# Populate the new row 81 with the new weekly record
$ws.Cells.Item(81, 1).Value = 10
$ws.Cells.Item(81, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(81, 3).Value = "La Araucanía"
$ws.Cells.Item(81, 4).Value = 44523
$ws.Cells.Item(81, 5).Value = 9
$ws.Cells.Item(81, 6).Value = "Fruta"
$ws.Cells.Item(81, 7).Value = 100107
$ws.Cells.Item(81, 8).Value = "Otros"
$ws.Cells.Item(81, 9).Value = 100107002
$ws.Cells.Item(81, 10).Value = "Chirimoya"
$ws.Cells.Item(81, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(81, 12).Value = "Primera"
$ws.Cells.Item(81, 13).Value = 95
$ws.Cells.Item(81, 14).Value = 3000
$ws.Cells.Item(81, 15).Value = 3000
$ws.Cells.Item(81, 16).Value = 3000
$ws.Cells.Item(81, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(81, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(81, 19).Value = 3000
$ws.Cells.Item(81, 20).Value = 1
